# Business/system exception handling run produced 7 additional processed
# email rows (rows 17-23) appended to the Partnership_Emails report, with
# the same Sender/Company/Address/VAT/Email/Subject values as the existing
# "no-attachment" rows (5-16) — only the "Date Processed" timestamp differs
# per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row number, date-processed serial value) for each newly appended row
$newRows = @(
    @{ Row = 17; Date = 45855.999027777776 },
    @{ Row = 18; Date = 45856.032638888886 },
    @{ Row = 19; Date = 45856.059907407405 },
    @{ Row = 20; Date = 45856.071342592593 },
    @{ Row = 21; Date = 45856.097303240742 },
    @{ Row = 22; Date = 45856.104791666665 },
    @{ Row = 23; Date = 45856.108865740738 }
)

foreach ($entry in $newRows) {
    $targetRow = $entry.Row
    $srcRow = $targetRow - 1

    # Duplicate the immediately preceding row (copy+insert) so the new row
    # inherits its formatting/styles (e.g. the date number format on column A)
    # exactly like the rest of the "no attachment" block below it.
    $ws.Rows($srcRow).Copy()
    $ws.Rows($targetRow).Insert()

    $ws.Cells.Item($targetRow, 1).Value = $entry.Date
    $ws.Cells.Item($targetRow, 2).Value = "Moris Mwai"
    $ws.Cells.Item($targetRow, 3).Value = "Tech-Neo GmbH"
    $ws.Cells.Item($targetRow, 4).Value = "Am main City, Germany"
    $ws.Cells.Item($targetRow, 5).Value = "DE1567890"
    $ws.Cells.Item($targetRow, 6).Value = "morismwai1@gmail.com"
    $ws.Cells.Item($targetRow, 7).Value = "Partnership Offer"
}
